$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.684.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "'2.674.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.64%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'513.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("D6").Value = "'142.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.93%  "
$ws.Range("D9").Value = "'2.666.90"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.06%  "
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("E11").Value = "  +3.94%  "
$ws.Range("D12").Value = "'0.334"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "'3.134.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.23%  "
$ws.Range("D15").Value = "'58.747.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "'20.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "'0.0000137"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("D18").Value = "'2.665.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.19%  "
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "'343.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.36%  "
$ws.Range("D21").Value = "'10.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.07%  "
$ws.Range("D22").Value = "'6.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.58%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("D25").Value = "'0.418"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.74%  "
$ws.Range("D26").Value = "'2.762.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.52%  "
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("D29").Value = "'0.0₃0804"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.67%  "
$ws.Range("D30").Value = "'7.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.26%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  +8.51%  "
$ws.Range("D33").Value = "'18.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("E34").Value = "  +2.02%  "
$ws.Range("D35").Value = "'149.63"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("E36").Value = "  +11.90%  "
$ws.Range("E37").Value = "  +1.41%  "
$ws.Range("E38").Value = "  +2.91%  "
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").Value = "'0.843"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.97%  "
$ws.Range("E41").Value = "  +4.02%  "
$ws.Range("D42").Value = "'1.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("E43").Value = "  +1.45%  "
$ws.Range("D44").Value = "'277.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.31%  "
$ws.Range("D45").Value = "'0.997"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").Value = "'19.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.98%  "
$ws.Range("D48").Value = "'0.0530"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("D49").Value = "'0.0230"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "'1.990.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.23%  "
